$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Abril de 2020 a las 09:22"

# --- Countries list reorder: Afganistan now comes right before Cuba ---
# (row 82 previously held Cuba, row 83 previously held Afganistan)
$ws.Range("A82").Value = "Afganistan"
$ws.Range("A83").Value = "Cuba"

# --- Hungria (row 62) updated stats ---
$ws.Range("E62").Value = 1477
$ws.Range("G62").Value = 17
$ws.Range("H62").Value = 189

# --- Afganistan (now row 82) updated stats ---
$ws.Range("B82").Value = 993
$ws.Range("C82").Value = 60
$ws.Range("D82").Value = 131
$ws.Range("E82").Value = 830
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 2
$ws.Range("H82").Value = 32

# --- Cuba (now row 83) updated stats ---
$ws.Range("B83").Value = 986
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 227
$ws.Range("E83").Value = 727
$ws.Range("F83").Value = 15
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 32

# --- Letonia (row 90) updated stats ---
$ws.Range("B90").Value = 727
$ws.Range("C90").Value = 15
$ws.Range("E90").Value = 634

# --- Sri Lanka (row 119) updated stats ---
$ws.Range("D119").Value = 91
$ws.Range("E119").Value = 156

# --- Islas Feroe (row 125) updated stats ---
$ws.Range("B125").Value = 185
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = 176
$ws.Range("E125").Value = 9
